# Расчет анкерных болтов — обновление количества болтов (E23): 36 -> 24
# Все остальные ячейки (H/I/J/L, AH19/AH25, E25, E38/E39, T41, D/E/F43:49, H55/J55/L55 и т.д.)
# являются формулами, зависящими от E23 (напрямую или через E25 = IF(E23=1,12,E23)),
# поэтому значение пересчитается автоматически движком при сохранении.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Фланец")

$ws.Range("E23").Value = 24
